$d = $word.ActiveDocument

# 1) Remove the old "_GoBack" bookmark that currently sits in the
#    "Finish division of pot for winners of the pool (Sean)" paragraph.
#    The diff moves this bookmark up into the "Worked on end of season..."
#    paragraph, so the old location must no longer carry it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Locate the paragraph that needs to be split into multiple runs and
#    gain the (new) "_GoBack" bookmark.
$target = "Worked on end of season and division of pot at the end of the season (Sean)"
$found = $false
foreach ($p in $d.Paragraphs) {
    $pr = $p.Range
    # Trim the trailing paragraph mark before comparing.
    $txt = $pr.Text
    if ($txt.Length -gt 0) {
        $txt = $txt.TrimEnd([char]13, [char]7)
    }
    if ($txt -eq $target) {
        $start = $pr.Start
        $end = $start + $target.Length
        $r = $d.Range($start, $end)

        $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Worked on end of season </w:t></w:r><w:r><w:t>division</w:t></w:r><w:r><w:t xml:space="preserve"> of pot for winners</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> (Sean)</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@
        $r.InsertXML($xml)
        $found = $true
        break
    }
}

if (-not $found) {
    throw "target paragraph not found"
}
